$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"

$ws.Range("D1").Select() | Out-Null
